$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.084.21'
$ws.Range('E2').Value = '  +2.53%  '
$ws.Range('D3').Value = '2.367.65'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.630'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.613'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0915'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.43'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.11%  '
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.974'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.85%  '
$ws.Range('D15').Value = '2.728.51'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.05%  '
$ws.Range('D17').Value = '2.380.46'
$ws.Range('E17').Value = '  +1.82%  '
$ws.Range('D18').Value = '45.080.21'
$ws.Range('E18').Value = '  +2.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.25'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.46%  '
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.22'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '258.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.44%  '
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('E30').Value = '  +8.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '37.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '168.60'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('E34').Value = '  +6.15%  '
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('E36').Value = '  +2.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.71'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.91'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.30%  '
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('E40').Value = '  -3.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.75'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.82'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.71%  '
$ws.Range('D43').Value = '1.896.02'
$ws.Range('E43').Value = '  +13.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '69.44'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.228'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.31%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('B47').Value = 'Celestia'
$ws.Range('C47').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '80.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.61'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '111.77'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.17'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.94%  '
